$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = "Handed back: in sync with en-US"
$overview.Range("F7").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = "Handed back: in sync with en-US"
$zhcn.Range("I7").Value = "7d77089f-af67-4d13-bf4f-e2576eac4631.md"
$zhcn.Range("J7").Value = "7d77089f-af67-4d13-bf4f-e2576eac4631.618a3581a5667c463aae1b1be488427e084e03bc.zh-cn.xlf"
$zhcn.Range("K7").Value = "2016-11-29 03:07:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = "Handed back: in sync with en-US"
$dede.Range("I7").Value = "7d77089f-af67-4d13-bf4f-e2576eac4631.md"
$dede.Range("J7").Value = "7d77089f-af67-4d13-bf4f-e2576eac4631.618a3581a5667c463aae1b1be488427e084e03bc.de-de.xlf"
$dede.Range("K7").Value = "2016-11-29 03:07:30"
